# Updated symbol list on Fri Dec 16 21:37:34 UTC 2022 with GitHub Actions
#
# All data cells in column B/C/D/E on this sheet are stored as text
# (inlineStr) in the source workbook, so every cell we touch gets its
# NumberFormat forced to "@" (Text) before the value is assigned. This
# keeps Excel from re-interpreting numeric-looking strings (prices,
# rank-prefixed labels, etc.) as actual numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every cell we are about to touch to Text format so values stay
# as literal strings (matches the source t="inlineStr" cells).
$ws.Range("B10:E18").NumberFormat = "@"
$ws.Range("B42:E43").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# --- Simple price / label refreshes (same coin, updated numbers) ---
$ws.Range("D2").Value = "242.80"
$ws.Range("D4").Value = "5.739"
$ws.Range("D5").Value = "0.05805"
$ws.Range("D7").Value = "6.474"
$ws.Range("D8").Value = "1.321"
$ws.Range("D9").Value = "0.8009"

# --- Rows 10-18: the coin ranking list shifted up by one slot, with
#     "One" wrapping from rank 9 (row 10) down to rank 17 (row 18). ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1463"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07654"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03242"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03003"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09232"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001682"
$ws.Range("E15").Value = "14BitForexTokenBF"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.273"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04756"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005995"
$ws.Range("E18").Value = "17OneONE"

# --- More simple price refreshes ---
$ws.Range("D19").Value = "0.006206"
$ws.Range("D20").Value = "0.005387"
$ws.Range("D21").Value = "0.001064"
$ws.Range("D22").Value = "0.0001502"
$ws.Range("D23").Value = "3.694"

$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"

$ws.Range("D40").Value = "0.04295"
$ws.Range("D41").Value = "0.007068"

# --- Rows 42-43: CEJI and BKEXToken swap ranking positions. ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1056"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.003316"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Final batch of price refreshes ---
$ws.Range("D46").Value = "0.00005630"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("D48").Value = "0.7859"
$ws.Range("D49").Value = "0.09982"
$ws.Range("D50").Value = "0.00002102"
